$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price column D, Volume(1h) column E)
# Each value is written with a leading apostrophe so Excel stores it as
# literal text (matching the original inlineStr cell type) instead of
# auto-converting the numeric-looking / percent-looking text into a number.
$ws.Range("D2").Value = "'261.53"
$ws.Range("E2").Value = "'0.21%"
$ws.Range("D3").Value = "'26.73"
$ws.Range("E3").Value = "'-2.38%"
$ws.Range("D4").Value = "'4.702"
$ws.Range("E4").Value = "'0.01%"
$ws.Range("E5").Value = "'1.63%"
$ws.Range("D6").Value = "'6.734"
$ws.Range("E6").Value = "'0.94%"
$ws.Range("D7").Value = "'0.8503"
$ws.Range("E7").Value = "'0.58%"
$ws.Range("D8").Value = "'0.9122"
$ws.Range("E8").Value = "'-1.71%"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'0.11%"
$ws.Range("D10").Value = "'0.05090"
$ws.Range("E10").Value = "'4.83%"
$ws.Range("D11").Value = "'0.07094"
$ws.Range("E11").Value = "'-0.13%"
$ws.Range("D12").Value = "'0.03105"
$ws.Range("E12").Value = "'0.88%"
$ws.Range("D13").Value = "'0.09045"
$ws.Range("E13").Value = "'-0.24%"
$ws.Range("D14").Value = "'0.001538"
$ws.Range("E14").Value = "'0.60%"
$ws.Range("D15").Value = "'0.0006182"
$ws.Range("E15").Value = "'1.79%"
$ws.Range("D16").Value = "'0.006008"
$ws.Range("E16").Value = "'-0.38%"
$ws.Range("D17").Value = "'3.448"
$ws.Range("E17").Value = "'-0.06%"
$ws.Range("D18").Value = "'3.173"
$ws.Range("E18").Value = "'0.79%"
$ws.Range("E21").Value = "'0.39%"
$ws.Range("D22").Value = "'4.126"
$ws.Range("E22").Value = "'0.65%"
$ws.Range("D23").Value = "'0.04249"
$ws.Range("E23").Value = "'-0.11%"
$ws.Range("E24").Value = "'-3.35%"
$ws.Range("D25").Value = "'0.004056"
$ws.Range("E25").Value = "'6.92%"
$ws.Range("E26").Value = "'0.05%"
$ws.Range("E27").Value = "'4.12%"
$ws.Range("D40").Value = "'0.03966"
$ws.Range("E40").Value = "'2.42%"
$ws.Range("E41").Value = "'0.03%"
$ws.Range("D42").Value = "'0.004138"
$ws.Range("E42").Value = "'1.36%"
$ws.Range("E44").Value = "'-18.80%"
$ws.Range("E45").Value = "'0.19%"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("D48").Value = "'0.2483"
$ws.Range("E48").Value = "'83.26%"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("E50").Value = "'0.07%"
